$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 4 (week 3): add prepare link in column D
$ws.Range("D4").Value = "prep/p03.html"

# Row 3 (week 2): add slides link in column E
$ws.Range("E3").Value = "slides/slides.html#/session-02-the-mobile-connectivity-paradox-digital-wellbeing-as-a-dynamic-construct"

# Update active selection to E4
$ws.Range("E4").Select()
